# Applies the "Added comments to understand parameters" change:
#  - Adds/edits cell comments across the 4 sheets explaining the parameters
#  - Sets channel!D2 (Q_ini_value) to 0
#  - Leaves the "channel" sheet as the active sheet/tab with N1 selected

$wb = $excel.ActiveWorkbook

$wsGeneral   = $wb.Worksheets.Item(1)   # "general"
$wsPeatHydro = $wb.Worksheets.Item(2)   # "peat_hydro_prop"
$wsChannel   = $wb.Worksheets.Item(3)   # "channel"
$wsPeat      = $wb.Worksheets.Item(4)   # "peat"

# ---------------------------------------------------------------------
# Sheet "general" (comments1.xml)
# ---------------------------------------------------------------------
$wsGeneral.Range("A1").AddComment("Outermost loop. Number of days to simulate. ")

$cGeneralC1 = $wsGeneral.Range("C1").Comment
$null = $cGeneralC1.Text("Used as a fallback value   when normal_timestep doesnt converge.`nIt serves as a very rough but performant adaptive timestep solution to the convergence problems of the hydrology model .")

# ---------------------------------------------------------------------
# Sheet "peat_hydro_prop" (new comments2.xml)
# ---------------------------------------------------------------------
$wsPeatHydro.Range("A1").AddComment("This is the number to reference from the code (and not the excel column numbers!)")
$wsPeatHydro.Range("B1").AddComment("See parameterizations for the specific yield and transmissivity functions")

# ---------------------------------------------------------------------
# Sheet "channel" (new comments3.xml)
# ---------------------------------------------------------------------
$wsChannel.Range("A1").AddComment("Metres.`nA value of 0 means that the highest point of the block (aka ‘block head’) is level with the adjacent peat.`nA value of -0.3 would mean that  the block head is 30 cm below the adjacent peat.")
$wsChannel.Range("B1").AddComment("SI units.`nUsed to compute the amount of water flow through the blocks.  See paper.")
$wsChannel.Range("C1").AddComment("Metres.`nInitial value for the CWL at canals, relative to the peat surface height (DEM).`n")
$wsChannel.Range("D1").AddComment("m^3/s`nInitial value for the discharge. Only used when solving with the Preissmann scheme.")
$wsChannel.Range("E1").AddComment("Metres.`nDepth of channel bed.`nWe don’t have measurements of the slope of the channel bed and so we assume that the channel bed is parallel to the peat surface.`nThe CWL can never go below this value, or otherwise a numerical error occurs. Therefore, this value should always be large. 8 metres has always been fine. ")
$wsChannel.Range("F1").AddComment("Metres.`nCWL boundary conditions at the  downstream and/or upstream  nodes of each channel reach.")
$wsChannel.Range("G1").AddComment("m³/s`nDirichlet BC for the water discharge. Only used for the Preissmann scheme")
$wsChannel.Range("H1").AddComment("Metres.`nWidth of the channels.")
$wsChannel.Range("I1").AddComment("Metres.`nDistance between nodes of the channel network.")
$wsChannel.Range("J1").AddComment("Parameters used to describe the Manning friction term in the open channel flow equations")
$wsChannel.Range("N1").AddComment("Maximum number of iterations to avoid infinite loops in the Newton method computation step")
$wsChannel.Range("O1").AddComment("Analogous quantity for the inexact fo accelerated Newton method. `nThe inexact newton method needs less computations of the jacobian, and, for that reason, it is usually more efficient")

# channel!D2 (Q_ini_value) gets an explicit initial value of 0
$wsChannel.Range("D2").Value = 0

# ---------------------------------------------------------------------
# Sheet "peat" (comments4.xml)
# ---------------------------------------------------------------------
$cPeatC1 = $wsPeat.Range("C1").Comment
$null = $cPeatC1.Text("Meters.`nUsed only if boundary_condition=’dirichlet’.`nElse, if boundary_condition=’neumann’, this value will be ignored and Neumann boundary conditions applied`n")

$cPeatD1 = $wsPeat.Range("D1").Comment
$null = $cPeatD1.Text("Metres.`nOnly used if rectangular grid hydrology`n")

# ---------------------------------------------------------------------
# Selection / active-sheet state
# ---------------------------------------------------------------------
$null = $wsPeatHydro.Range("B1").Select()
$null = $wsPeat.Range("D1").Select()
$null = $wsChannel.Activate()
$null = $wsChannel.Range("N1").Select()
